{"js": "const replacements = [\n  [\"2024-12-07 Saturday\", \"2024-12-08 Sunday\"],\n  [\"78\u00d759=4602\", \"52\u00d766=3432\"],\n  [\"57\u00d712=684\", \"44\u00d737=1628\"],\n  [\"87\u00d734=2958\", \"68\u00d768=4624\"],\n  [\"52\u00d759=3068\", \"84\u00d720=1680\"],\n  [\"16\u00d733=528\", \"13\u00d773=949\"],\n  [\"54\u00d749=2646\", \"26\u00d768=1768\"],\n  [\"96\u00d754=5184\", \"71\u00d760=4260\"],\n  [\"28\u00d757=1596\", \"52\u00d715=780\"],\n  [\"41\u00d794=3854\", \"77\u00d762=4774\"],\n  [\"83\u00d743=3569\", \"30\u00d713=390\"],\n  [\"25\u00d743=1075\", \"41\u00d756=2296\"],\n  [\"17\u00d795=1615\", \"64\u00d760=3840\"],\n  [\"24\u00d735=840\", \"49\u00d757=2793\"],\n  [\"96\u00d798=9408\", \"55\u00d778=4290\"],\n  [\"67\u00d716=1072\", \"72\u00d735=2520\"],\n  [\"93\u00d773=6789\", \"34\u00d747=1598\"],\n  [\"99\u00d775=7425\", \"27\u00d727=729\"],\n  [\"83\u00d726=2158\", \"63\u00d766=4158\"],\n  [\"39\u00d763=2457\", \"46\u00d749=2254\"],\n  [\"51\u00d712=612\", \"37\u00d748=1776\"],\n  [\"73\u00d728=2044\", \"25\u00d717=425\"],\n  [\"60\u00d737=2220\", \"44\u00d763=2772\"],\n  [\"68\u00d739=2652\", \"97\u00d713=1261\"],\n  [\"81\u00d799=8019\", \"54\u00d782=4428\"],\n  [\"79\u00d711=869\", \"17\u00d799=1683\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-12-07 Saturday\"; New = \"2024-12-08 Sunday\" },\n    @{ Old = \"78\u00d759=4602\"; New = \"52\u00d766=3432\" },\n    @{ Old = \"57\u00d712=684\"; New = \"44\u00d737=1628\" },\n    @{ Old = \"87\u00d734=2958\"; New = \"68\u00d768=4624\" },\n    @{ Old = \"52\u00d759=3068\"; New = \"84\u00d720=1680\" },\n    @{ Old = \"16\u00d733=528\"; New = \"13\u00d773=949\" },\n    @{ Old = \"54\u00d749=2646\"; New = \"26\u00d768=1768\" },\n    @{ Old = \"96\u00d754=5184\"; New = \"71\u00d760=4260\" },\n    @{ Old = \"28\u00d757=1596\"; New = \"52\u00d715=780\" },\n    @{ Old = \"41\u00d794=3854\"; New = \"77\u00d762=4774\" },\n    @{ Old = \"83\u00d743=3569\"; New = \"30\u00d713=390\" },\n    @{ Old = \"25\u00d743=1075\"; New = \"41\u00d756=2296\" },\n    @{ Old = \"17\u00d795=1615\"; New = \"64\u00d760=3840\" },\n    @{ Old = \"24\u00d735=840\"; New = \"49\u00d757=2793\" },\n    @{ Old = \"96\u00d798=9408\"; New = \"55\u00d778=4290\" },\n    @{ Old = \"67\u00d716=1072\"; New = \"72\u00d735=2520\" },\n    @{ Old = \"93\u00d773=6789\"; New = \"34\u00d747=1598\" },\n    @{ Old = \"99\u00d775=7425\"; New = \"27\u00d727=729\" },\n    @{ Old = \"83\u00d726=2158\"; New = \"63\u00d766=4158\" },\n    @{ Old = \"39\u00d763=2457\"; New = \"46\u00d749=2254\" },\n    @{ Old = \"51\u00d712=612\"; New = \"37\u00d748=1776\" },\n    @{ Old = \"73\u00d728=2044\"; New = \"25\u00d717=425\" },\n    @{ Old = \"60\u00d737=2220\"; New = \"44\u00d763=2772\" },\n    @{ Old = \"68\u00d739=2652\"; New = \"97\u00d713=1261\" },\n    @{ Old = \"81\u00d799=8019\"; New = \"54\u00d782=4428\" },\n    @{ Old = \"79\u00d711=869\"; New = \"17\u00d799=1683\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
